$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45205 -> 45206) for every data row (rows 2 through 362).
$ws.Range("C2:C362").Value = 45206
